$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (locked cells); unprotect so the values/text can be
# updated, then restore protection at the end.
$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure footnote text
# (shared string used by cell A58): 2021-04-06 -> 2021-04-08
$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Refresh the Weight (column D) and Percent Change (column E) figures for
# each holding row (rows 2-55)
$ws.Range("D2").Value = 0.01610730416491809
$ws.Range("E2").Value = -0.001279590531030217
$ws.Range("D3").Value = 0.05215331805465782
$ws.Range("E3").Value = 0.006071251055836502
$ws.Range("D4").Value = 0.01483401759127845
$ws.Range("E4").Value = 0.01099961404862992
$ws.Range("D5").Value = 0.009581456021263174
$ws.Range("E5").Value = -0.002522905324658042
$ws.Range("D6").Value = 0.01550876532741932
$ws.Range("E6").Value = -0.006066522557701681
$ws.Range("D7").Value = 0.02044147853151432
$ws.Range("E7").Value = 0.003802837501828193
$ws.Range("D8").Value = 0.004246101561779876
$ws.Range("E8").Value = -0.01253211682659539
$ws.Range("D9").Value = 0.006479792014592415
$ws.Range("E9").Value = -0.01042587029510522
$ws.Range("D10").Value = 0.01394403378789588
$ws.Range("E10").Value = 0
$ws.Range("D11").Value = 0.00871122419035838
$ws.Range("E11").Value = 0.004454505622900395
$ws.Range("D12").Value = 0.01465653609944918
$ws.Range("E12").Value = -0.006770833333333282
$ws.Range("D13").Value = 0.003228381975030581
$ws.Range("E13").Value = -0.01517241379310352
$ws.Range("D14").Value = 0.006129790879092547
$ws.Range("E14").Value = 0.001141552511415345
$ws.Range("D15").Value = 0.01429066341118537
$ws.Range("E15").Value = 0.00122636029174461
$ws.Range("D16").Value = 0.01040413680731407
$ws.Range("E16").Value = -0.004952538175815091
$ws.Range("D17").Value = 0.0213049355312794
$ws.Range("E17").Value = 0.003627813234799993
$ws.Range("D18").Value = 0.008290484911878532
$ws.Range("E18").Value = -0.00241701579116993
$ws.Range("D19").Value = 0.01665247132799223
$ws.Range("E19").Value = -0.00391174133610428
$ws.Range("D20").Value = 0.01159800200369607
$ws.Range("E20").Value = 0.001590609916630292
$ws.Range("D21").Value = 0.007279985450334477
$ws.Range("E21").Value = -0.03006789524733278
$ws.Range("D22").Value = 0.0132679819762774
$ws.Range("E22").Value = -0.009397234528124354
$ws.Range("D23").Value = 0.01888889715863582
$ws.Range("E23").Value = -0.007501704932939224
$ws.Range("D24").Value = 0.009556010646090518
$ws.Range("E24").Value = -0.005791505791505669
$ws.Range("D25").Value = 0.02097538611607455
$ws.Range("E25").Value = 0.01048730484150773
$ws.Range("D26").Value = 0.01148524718496224
$ws.Range("E26").Value = 0.006222739168362912
$ws.Range("D27").Value = 0.02304899335233213
$ws.Range("E27").Value = 0.01399556756428555
$ws.Range("D28").Value = 0.05695311098019466
$ws.Range("E28").Value = 0.01923377638780299
$ws.Range("D29").Value = 0.02127166570324115
$ws.Range("E29").Value = 0.006321968841724868
$ws.Range("D30").Value = 0.03219476093720151
$ws.Range("E30").Value = 0.01434499110847676
$ws.Range("D31").Value = 0.01640880005399509
$ws.Range("E31").Value = 0.006526572473642744
$ws.Range("D32").Value = 0.0136807695750158
$ws.Range("E32").Value = -0.001325205406838004
$ws.Range("D33").Value = 0.02151762706100382
$ws.Range("E33").Value = 0.01415495955725854
$ws.Range("D34").Value = 0.042729718779622
$ws.Range("E34").Value = 0.005091490511516028
$ws.Range("D35").Value = 0.01106873820010485
$ws.Range("E35").Value = -0.009655172413793212
$ws.Range("D36").Value = 0.009369877726702548
$ws.Range("E36").Value = 0.008459271932325763
$ws.Range("D37").Value = 0.01197522969063068
$ws.Range("E37").Value = -0.00776892430278886
$ws.Range("D38").Value = 0.007225214280275337
$ws.Range("E38").Value = 0.008055995773903923
$ws.Range("D39").Value = 0.0118468895796036
$ws.Range("E39").Value = -0.008739076154806291
$ws.Range("D40").Value = 0.01745431871312051
$ws.Range("E40").Value = 0.002704268881591698
$ws.Range("D41").Value = 0.0169671033920021
$ws.Range("E41").Value = -0.003419290497221783
$ws.Range("D42").Value = 0.03386365948133801
$ws.Range("E42").Value = 0.01851277860745593
$ws.Range("D43").Value = 0.01120423482289924
$ws.Range("E43").Value = 0.0006472491909386147
$ws.Range("D44").Value = 0.02162020372966859
$ws.Range("E44").Value = 0.006521639987230321
$ws.Range("D45").Value = 0.01381092266902392
$ws.Range("E45").Value = 0.006259586286946117
$ws.Range("D46").Value = 0.00798981599749465
$ws.Range("E46").Value = 0.01909641360037262
$ws.Range("D47").Value = 0.01319145501044565
$ws.Range("E47").Value = 0.0003038062584088674
$ws.Range("D48").Value = 0.009582282995956285
$ws.Range("E48").Value = -0.01389137838322274
$ws.Range("D49").Value = 0.015038216727039
$ws.Range("E49").Value = 0.01598984771573608
$ws.Range("D50").Value = 0.008303748313687277
$ws.Range("E50").Value = -0.006576805365631322
$ws.Range("D51").Value = 0.01114017609090208
$ws.Range("E51").Value = -0.007594662037539335
$ws.Range("D52").Value = 0.008537623118242943
$ws.Range("E52").Value = -0.004712728464879912
$ws.Range("D53").Value = 0.1381202954284394
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 0.04386814486484657
$ws.Range("E54").Value = 0.006655984222852185
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 0.003609585501836099

# Restore sheet protection
$ws.Protect()
